$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 3.272327238179451
$ws.Cells.Item(2, 3).Value = 1.626987699542094
$ws.Cells.Item(2, 4).Value = 0.7210945179870265
$ws.Cells.Item(2, 5).Value = 0.5333859586016987
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 6.15379541431027

$ws.Cells.Item(3, 2).Value = 3.272327238179451
$ws.Cells.Item(3, 3).Value = 1.626987699542094
$ws.Cells.Item(3, 4).Value = 0.7210945179870265
$ws.Cells.Item(3, 5).Value = 0.5333859586016987
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 6.15379541431027

$ws.Cells.Item(4, 2).Value = 3.272327238179451
$ws.Cells.Item(4, 3).Value = 1.626987699542094
$ws.Cells.Item(4, 4).Value = 0.7210945179870265
$ws.Cells.Item(4, 5).Value = 0.5333859586016987
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 6.15379541431027

$ws.Cells.Item(5, 2).Value = 1.445647641019636
$ws.Cells.Item(5, 3).Value = 1.626987699542094
$ws.Cells.Item(5, 4).Value = 0.7210945179870265
$ws.Cells.Item(5, 5).Value = 0.5333859586016987
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 4.327115817150455

$ws.Cells.Item(6, 2).Value = 0.6545652718822623
$ws.Cells.Item(6, 3).Value = 0.3048912486333797
$ws.Cells.Item(6, 4).Value = 0.1496068669990043
$ws.Cells.Item(6, 5).Value = 0.5333859586016987
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.642449346116345

$ws.Cells.Item(7, 2).Value = 0.1169995834814548
$ws.Cells.Item(7, 3).Value = 0.002658071450198252
$ws.Cells.Item(7, 4).Value = 3.223369029078222
$ws.Cells.Item(7, 5).Value = 0.5333859586016987
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 3.876412642611573

$ws.Cells.Item(8, 2).Value = 3.272327238179451
$ws.Cells.Item(8, 3).Value = 1.626987699542094
$ws.Cells.Item(8, 4).Value = 3.223369029078222
$ws.Cells.Item(8, 5).Value = 0.5333859586016987
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 8.656069925401464

$ws.Cells.Item(9, 2).Value = 0.1169995834814548
$ws.Cells.Item(9, 3).Value = 0.3048912486333797
$ws.Cells.Item(9, 4).Value = 0.7210945179870265
$ws.Cells.Item(9, 5).Value = 0.5333859586016987
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.67637130870356

$ws.Cells.Item(10, 2).Value = 3.272327238179451
$ws.Cells.Item(10, 3).Value = 1.626987699542094
$ws.Cells.Item(10, 4).Value = 0.7210945179870265
$ws.Cells.Item(10, 5).Value = 13.86384647080068
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 19.48425592650926

$ws.Cells.Item(11, 2).Value = 3.272327238179451
$ws.Cells.Item(11, 3).Value = 1.626987699542094
$ws.Cells.Item(11, 4).Value = 0.1496068669990043
$ws.Cells.Item(11, 5).Value = 0.5333859586016987
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 5.582307763322248

$ws.Cells.Item(12, 2).Value = 3.272327238179451
$ws.Cells.Item(12, 3).Value = 1.626987699542094
$ws.Cells.Item(12, 4).Value = 18.71679738969934
$ws.Cells.Item(12, 5).Value = 0.5333859586016987
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 24.14949828602258

$ws.Cells.Item(13, 2).Value = 3.272327238179451
$ws.Cells.Item(13, 3).Value = 1.626987699542094
$ws.Cells.Item(13, 4).Value = 0.1496068669990043
$ws.Cells.Item(13, 5).Value = 13.86384647080068
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 18.91276827552123

$ws.Cells.Item(14, 2).Value = 1.445647641019636
$ws.Cells.Item(14, 3).Value = 1.626987699542094
$ws.Cells.Item(14, 4).Value = 0.7210945179870265
$ws.Cells.Item(14, 5).Value = 0.5333859586016987
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 4.327115817150455

$ws.Cells.Item(15, 2).Value = 3.272327238179451
$ws.Cells.Item(15, 3).Value = 1.626987699542094
$ws.Cells.Item(15, 4).Value = 0.7210945179870265
$ws.Cells.Item(15, 5).Value = 0.5333859586016987
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 6.15379541431027

$ws.Cells.Item(16, 2).Value = 3.272327238179451
$ws.Cells.Item(16, 3).Value = 1.626987699542094
$ws.Cells.Item(16, 4).Value = 18.71679738969934
$ws.Cells.Item(16, 5).Value = 0.5333859586016987
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 24.14949828602258

$ws.Cells.Item(17, 2).Value = 3.272327238179451
$ws.Cells.Item(17, 3).Value = 1.626987699542094
$ws.Cells.Item(17, 4).Value = 0.7210945179870265
$ws.Cells.Item(17, 5).Value = 0.5333859586016987
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 6.15379541431027

$ws.Cells.Item(18, 2).Value = 1.445647641019636
$ws.Cells.Item(18, 3).Value = 0.04103571897497393
$ws.Cells.Item(18, 4).Value = 0.7210945179870265
$ws.Cells.Item(18, 5).Value = 0.5333859586016987
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 2.741163836583335

$ws.Cells.Item(19, 2).Value = 1.445647641019636
$ws.Cells.Item(19, 3).Value = 1.626987699542094
$ws.Cells.Item(19, 4).Value = 18.71679738969934
$ws.Cells.Item(19, 5).Value = 0.5333859586016987
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 22.32281868886277

$ws.Cells.Item(20, 2).Value = 1.445647641019636
$ws.Cells.Item(20, 3).Value = 1.626987699542094
$ws.Cells.Item(20, 4).Value = 3.223369029078222
$ws.Cells.Item(20, 5).Value = 0.5333859586016987
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 6.82939032824165

$ws.Cells.Item(21, 2).Value = 3.272327238179451
$ws.Cells.Item(21, 3).Value = 1.626987699542094
$ws.Cells.Item(21, 4).Value = 0.7210945179870265
$ws.Cells.Item(21, 5).Value = 0.5333859586016987
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 6.15379541431027

$ws.Cells.Item(22, 2).Value = 0.01253208636536152
$ws.Cells.Item(22, 3).Value = 0.3048912486333797
$ws.Cells.Item(22, 4).Value = 3.223369029078222
$ws.Cells.Item(22, 5).Value = 13.86384647080068
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 17.40463883487765

$ws.Cells.Item(23, 2).Value = 0.6545652718822623
$ws.Cells.Item(23, 3).Value = 9.98352242611593
$ws.Cells.Item(23, 4).Value = 3.223369029078222
$ws.Cells.Item(23, 5).Value = 13.86384647080068
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 27.7253031978771

$ws.Cells.Item(24, 2).Value = 3.272327238179451
$ws.Cells.Item(24, 3).Value = 1.626987699542094
$ws.Cells.Item(24, 4).Value = 0.1496068669990043
$ws.Cells.Item(24, 5).Value = 0.5333859586016987
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 5.582307763322248

